{"js": "// Replace the two-digit multiplication problems/answers throughout the\n// document body (including inside the table cells) with the updated\n// values, as described by the diff.\nconst replacements = [\n  [\"37\u00d773=2701\", \"42\u00d785=3570\"],\n  [\"74\u00d751=3774\", \"47\u00d784=3948\"],\n  [\"30\u00d732=960\", \"40\u00d751=2040\"],\n  [\"78\u00d733=2574\", \"50\u00d714=700\"],\n  [\"26\u00d744=1144\", \"48\u00d765=3120\"],\n  [\"83\u00d724=1992\", \"33\u00d712=396\"],\n  [\"98\u00d785=8330\", \"37\u00d762=2294\"],\n  [\"89\u00d750=4450\", \"80\u00d773=5840\"],\n  [\"11\u00d720=220\", \"75\u00d722=1650\"],\n  [\"80\u00d796=7680\", \"72\u00d740=2880\"],\n  [\"57\u00d756=3192\", \"25\u00d794=2350\"],\n  [\"34\u00d730=1020\", \"46\u00d730=1380\"],\n  [\"73\u00d758=4234\", \"63\u00d736=2268\"],\n  [\"34\u00d736=1224\", \"61\u00d772=4392\"],\n  [\"84\u00d728=2352\", \"88\u00d784=7392\"],\n  [\"12\u00d786=1032\", \"96\u00d773=7008\"],\n  [\"23\u00d730=690\", \"63\u00d768=4284\"],\n  [\"68\u00d739=2652\", \"47\u00d793=4371\"],\n  [\"65\u00d725=1625\", \"87\u00d732=2784\"],\n  [\"99\u00d770=6930\", \"53\u00d786=4558\"],\n  [\"18\u00d793=1674\", \"19\u00d762=1178\"],\n  [\"55\u00d786=4730\", \"66\u00d719=1254\"],\n  [\"63\u00d730=1890\", \"30\u00d792=2760\"],\n  [\"40\u00d787=3480\", \"15\u00d764=960\"],\n  [\"12\u00d791=1092\", \"93\u00d784=7812\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems/answers throughout the\n# document (including inside the table cells) with the updated values,\n# as described by the diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Before = \"37\u00d773=2701\"; After = \"42\u00d785=3570\" },\n    @{ Before = \"74\u00d751=3774\"; After = \"47\u00d784=3948\" },\n    @{ Before = \"30\u00d732=960\";  After = \"40\u00d751=2040\" },\n    @{ Before = \"78\u00d733=2574\"; After = \"50\u00d714=700\" },\n    @{ Before = \"26\u00d744=1144\"; After = \"48\u00d765=3120\" },\n    @{ Before = \"83\u00d724=1992\"; After = \"33\u00d712=396\" },\n    @{ Before = \"98\u00d785=8330\"; After = \"37\u00d762=2294\" },\n    @{ Before = \"89\u00d750=4450\"; After = \"80\u00d773=5840\" },\n    @{ Before = \"11\u00d720=220\";  After = \"75\u00d722=1650\" },\n    @{ Before = \"80\u00d796=7680\"; After = \"72\u00d740=2880\" },\n    @{ Before = \"57\u00d756=3192\"; After = \"25\u00d794=2350\" },\n    @{ Before = \"34\u00d730=1020\"; After = \"46\u00d730=1380\" },\n    @{ Before = \"73\u00d758=4234\"; After = \"63\u00d736=2268\" },\n    @{ Before = \"34\u00d736=1224\"; After = \"61\u00d772=4392\" },\n    @{ Before = \"84\u00d728=2352\"; After = \"88\u00d784=7392\" },\n    @{ Before = \"12\u00d786=1032\"; After = \"96\u00d773=7008\" },\n    @{ Before = \"23\u00d730=690\";  After = \"63\u00d768=4284\" },\n    @{ Before = \"68\u00d739=2652\"; After = \"47\u00d793=4371\" },\n    @{ Before = \"65\u00d725=1625\"; After = \"87\u00d732=2784\" },\n    @{ Before = \"99\u00d770=6930\"; After = \"53\u00d786=4558\" },\n    @{ Before = \"18\u00d793=1674\"; After = \"19\u00d762=1178\" },\n    @{ Before = \"55\u00d786=4730\"; After = \"66\u00d719=1254\" },\n    @{ Before = \"63\u00d730=1890\"; After = \"30\u00d792=2760\" },\n    @{ Before = \"40\u00d787=3480\"; After = \"15\u00d764=960\" },\n    @{ Before = \"12\u00d791=1092\"; After = \"93\u00d784=7812\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $r.Before,      # FindText\n        $true,          # MatchCase\n        $true,          # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap (wdFindContinue)\n        $false,         # Format\n        $r.After,       # ReplaceWith\n        2               # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
